$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/attribution-code-vs"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Include from Attribution Code sheet ---
$codeSheet = $wb.Worksheets.Item("Include from Attribution Code")
$codeSheet.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/attribution-code-system"

# --- Include from Attribution Sour(ce) sheet ---
$sourceSheet = $wb.Worksheets.Item("Include from Attribution Sour")
$sourceSheet.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/attribution-source"

# --- Include from Attribution Prod(uct) sheet ---
$productSheet = $wb.Worksheets.Item("Include from Attribution Prod")
$productSheet.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/attribution-product"
